$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.192899
$ws.Range("H2").Value = 0.578697
$ws.Range("I2").Value = 0.04416426199014034
$ws.Range("J2").Value = 0.04416426199014033
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.051093
$ws.Range("N2").Value = 0.153279
$ws.Range("O2").Value = 0.01450579975525089
$ws.Range("P2").Value = 0.01450579975525089
$ws.Range("Q2").Value = 0.009855788607000001
$ws.Range("R2").Value = 0.08870209746300001
$ws.Range("S2").Value = 0.0006406379407674137
$ws.Range("T2").Value = 0.0006406379407674138
$ws.Range("G3").Value = 0.192899
$ws.Range("H3").Value = 0.578697
$ws.Range("I3").Value = 0.04416426199014034
$ws.Range("J3").Value = 0.04416426199014033
$ws.Range("O3").Value = 0.2313022967634575
$ws.Range("P3").Value = 0.2313022967634575
$ws.Range("Q3").Value = 0.1571555225963333
$ws.Range("R3").Value = 1.414399703367
$ws.Range("S3").Value = 0.01021529523318252
$ws.Range("T3").Value = 0.01021529523318253
$ws.Range("G4").Value = 0.192899
$ws.Range("H4").Value = 0.578697
$ws.Range("I4").Value = 0.04416426199014034
$ws.Range("J4").Value = 0.04416426199014033
$ws.Range("M4").Value = 2.656449666666667
$ws.Range("N4").Value = 7.969348999999999
$ws.Range("O4").Value = 0.7541919034812916
$ws.Range("P4").Value = 0.7541919034812917
$ws.Range("Q4").Value = 0.5124264842503333
$ws.Range("R4").Value = 4.611838358252999
$ws.Range("S4").Value = 0.0333083288161904
$ws.Range("T4").Value = 0.0333083288161904
$ws.Range("I5").Value = 0.7315465959037607
$ws.Range("J5").Value = 0.7315465959037606
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.051093
$ws.Range("N5").Value = 0.153279
$ws.Range("O5").Value = 0.01450579975525089
$ws.Range("P5").Value = 0.01450579975525089
$ws.Range("Q5").Value = 0.163253460615
$ws.Range("R5").Value = 1.469281145535
$ws.Range("S5").Value = 0.01061166843181539
$ws.Range("T5").Value = 0.01061166843181539
$ws.Range("I6").Value = 0.7315465959037607
$ws.Range("J6").Value = 0.7315465959037606
$ws.Range("O6").Value = 0.2313022967634575
$ws.Range("P6").Value = 0.2313022967634575
$ws.Range("S6").Value = 0.1692084078220288
$ws.Range("T6").Value = 0.1692084078220288
$ws.Range("I7").Value = 0.7315465959037607
$ws.Range("J7").Value = 0.7315465959037606
$ws.Range("M7").Value = 2.656449666666667
$ws.Range("N7").Value = 7.969348999999999
$ws.Range("O7").Value = 0.7541919034812916
$ws.Range("P7").Value = 0.7541919034812917
$ws.Range("Q7").Value = 8.487945531342778
$ws.Range("R7").Value = 76.391509782085
$ws.Range("S7").Value = 0.5517265196499165
$ws.Range("T7").Value = 0.5517265196499165
$ws.Range("G8").Value = 0.9796416666666667
$ws.Range("H8").Value = 2.938925
$ws.Range("I8").Value = 0.224289142106099
$ws.Range("J8").Value = 0.224289142106099
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.051093
$ws.Range("N8").Value = 0.153279
$ws.Range("O8").Value = 0.01450579975525089
$ws.Range("P8").Value = 0.01450579975525089
$ws.Range("Q8").Value = 0.05005283167500001
$ws.Range("R8").Value = 0.450475485075
$ws.Range("S8").Value = 0.003253493382668083
$ws.Range("T8").Value = 0.003253493382668083
$ws.Range("G9").Value = 0.9796416666666667
$ws.Range("H9").Value = 2.938925
$ws.Range("I9").Value = 0.224289142106099
$ws.Range("J9").Value = 0.224289142106099
$ws.Range("O9").Value = 0.2313022967634575
$ws.Range("P9").Value = 0.2313022967634575
$ws.Range("Q9").Value = 0.7981176578527778
$ws.Range("R9").Value = 7.183058920675
$ws.Range("S9").Value = 0.0518785937082462
$ws.Range("T9").Value = 0.05187859370824621
$ws.Range("G10").Value = 0.9796416666666667
$ws.Range("H10").Value = 2.938925
$ws.Range("I10").Value = 0.224289142106099
$ws.Range("J10").Value = 0.224289142106099
$ws.Range("M10").Value = 2.656449666666667
$ws.Range("N10").Value = 7.969348999999999
$ws.Range("O10").Value = 0.7541919034812916
$ws.Range("P10").Value = 0.7541919034812917
$ws.Range("Q10").Value = 2.602368778869445
$ws.Range("R10").Value = 23.421319009825
$ws.Range("S10").Value = 0.1691570550151847
$ws.Range("T10").Value = 0.1691570550151847
